$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.499.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.962.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4766"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4074"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08474"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.062"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.966.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.645"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001072"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06597"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.826"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.513.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.291"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.192.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.949"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.166"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9894"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09602"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.456"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.621"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.660"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02356"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.916"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06237"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.254"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6220"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.006"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1927"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.349"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5986"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.076"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.400"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000331"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06843"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.70%  "
